# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.85 = 6710.77 pesos`n✅ 6710.77 pesos = 1.85 = 954.72 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the N10/O10 and N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 539.998
$wsTasas.Range("O10").Value = 3623.8
$wsTasas.Range("N12").Value = 3633.99
$wsTasas.Range("O12").Value = 516.999
